# Updated cryptos list on Tue Nov 12 11:23:45 UTC 2024 with GitHub Actions
#
# Refreshes the "Price" and "Volume(1h)" columns for every coin row with the
# latest scraped figures. A couple of neighbouring coins (rows 40/41 and
# 48/49) also swapped rank order, so their Coin name / Link cells are updated
# too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, column index, new cell text (exactly as it should appear).
$updates = @(
    @(2, 4, "87.222.84"),
    @(2, 5, "  +6.28%  "),
    @(3, 4, "3.280.84"),
    @(3, 5, "  +2.87%  "),
    @(4, 4, "0.999"),
    @(4, 5, "  -0.09%  "),
    @(5, 4, "214.36"),
    @(5, 5, "  -0.62%  "),
    @(6, 4, "629.24"),
    @(6, 5, "  +0.79%  "),
    @(7, 4, "0.392"),
    @(7, 5, "  +36.50%  "),
    @(8, 5, "  -0.06%  "),
    @(9, 4, "0.643"),
    @(9, 5, "  +9.80%  "),
    @(10, 4, "3.273.40"),
    @(10, 5, "  +2.78%  "),
    @(11, 4, "0.587"),
    @(11, 5, "  -0.79%  "),
    @(12, 4, "0.0000267"),
    @(12, 5, "  +3.16%  "),
    @(13, 5, "  +5.67%  "),
    @(14, 4, "34.73"),
    @(14, 5, "  +9.43%  "),
    @(15, 4, "3.885.85"),
    @(15, 5, "  +3.03%  "),
    @(16, 4, "5.27"),
    @(16, 5, "  -1.03%  "),
    @(17, 4, "86.966.04"),
    @(17, 5, "  +6.55%  "),
    @(18, 4, "3.281.51"),
    @(18, 5, "  +3.11%  "),
    @(19, 4, "14.23"),
    @(19, 5, "  +1.13%  "),
    @(20, 5, "  -7.28%  "),
    @(21, 4, "9.25"),
    @(21, 5, "  +3.10%  "),
    @(22, 4, "438.35"),
    @(22, 5, "  +0.63%  "),
    @(23, 5, "  +4.54%  "),
    @(24, 4, "7.27"),
    @(24, 5, "  +0.01%  "),
    @(25, 4, "5.24"),
    @(25, 5, "  -1.88%  "),
    @(26, 4, "12.27"),
    @(26, 5, "  +11.08%  "),
    @(27, 4, "3.473.26"),
    @(27, 5, "  +3.51%  "),
    @(28, 4, "77.12"),
    @(28, 5, "  +0.56%  "),
    @(29, 4, "0.0000132"),
    @(29, 5, "  +6.85%  "),
    @(30, 4, "0.999"),
    @(30, 5, "  -0.01%  "),
    @(31, 4, "0.179"),
    @(31, 5, "  +28.26%  "),
    @(32, 4, "0.996"),
    @(32, 5, "  -0.43%  "),
    @(33, 4, "8.98"),
    @(33, 5, "  -1.24%  "),
    @(34, 4, "556.51"),
    @(34, 5, "  -5.24%  "),
    @(35, 5, "  -3.77%  "),
    @(36, 4, "1.99"),
    @(36, 5, "  -1.01%  "),
    @(37, 4, "6.95"),
    @(37, 5, "  +12.72%  "),
    @(38, 5, "  -10.75%  "),
    @(39, 4, "22.70"),
    @(39, 5, "  -0.51%  "),
    @(40, 2, "WhiteBITCoin"),
    @(40, 3, "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"),
    @(40, 4, "21.73"),
    @(40, 5, "  +4.51%  "),
    @(41, 2, "FirstDigitalUSD"),
    @(41, 3, "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"),
    @(41, 4, "0.995"),
    @(41, 5, "  -0.35%  "),
    @(42, 4, "0.402"),
    @(42, 5, "  -1.74%  "),
    @(43, 4, "2.03"),
    @(43, 5, "  -1.13%  "),
    @(44, 4, "2.99"),
    @(44, 5, "  -2.28%  "),
    @(45, 5, "  -0.06%  "),
    @(46, 4, "152.88"),
    @(46, 5, "  -5.18%  "),
    @(47, 4, "181.38"),
    @(47, 5, "  -3.48%  "),
    @(48, 2, "ImmutableX"),
    @(48, 3, "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"),
    @(48, 4, "1.36"),
    @(48, 5, "  +1.54%  "),
    @(49, 2, "OKB"),
    @(49, 3, "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"),
    @(49, 4, "45.13"),
    @(49, 5, "  +0.87%  "),
    @(50, 4, "4.27"),
    @(50, 5, "  +1.29%  "),
    @(51, 4, "0.750"),
    @(51, 5, "  -3.14%  "),
)

foreach ($u in $updates) {
    $r = $u[0]
    $col = $u[1]
    $text = $u[2]
    $cell = $ws.Cells.Item($r, $col)

    # Numeric-looking text (e.g. "0.999", "214.36") would otherwise be
    # auto-coerced into a real number by Excel on assignment, which would
    # change the cell from text to numeric and could drop things like
    # trailing zeros. Force a text number format while writing the value,
    # then clear the format again so no extra styling is left behind -
    # only the underlying string representation is preserved.
    $isNumeric = $text -match '^[0-9]+(\.[0-9]+)?$'
    if ($isNumeric) {
        $cell.NumberFormat = "@"
        $cell.Value2 = $text
        $cell.ClearFormats()
    } else {
        $cell.Value2 = $text
    }
}
